$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165100455284119
$ws.Range("B1").Value = 1.821740746498108
$ws.Range("C1").Value = 4.352670192718506
$ws.Range("D1").Value = 2.280625343322754
$ws.Range("E1").Value = 0.3693938851356506
